$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: torta -> Materias primas (C2)
$ws.Range("C2").Value = "huevos,vainilla,leche,harina,"

# Row 3: kuchen -> Materias primas (C3)
$ws.Range("C3").Value = "harina,manzana,huevos,"

# Row 4: queque -> Materias primas (C4)
$ws.Range("C4").Value = "huevos,vainilla,harina,"

# Row 6: pie de limon -> Materias primas (C6)
$ws.Range("C6").Value = "crema,harina, limon, huevos,merengue,"
